$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classroom")

# Update rubric point values ("save and restore" scoring adjustments)
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 6
$ws.Range("E8").Value = 8
$ws.Range("E9").Value = 8
$ws.Range("E18").Value = 8
$ws.Range("E19").Value = 8

# Move the active selection from E22 to E6 and scroll the view back to the top
$ws.Activate() | Out-Null
$ws.Range("E6").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Reposition the application window (mirrors the saved workbookView xWindow/yWindow)
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 500
